$d = $word.ActiveDocument

$replacements = @(
    @("89÷3=", "65÷7="),
    @("24÷8=", "70÷2="),
    @("88÷8=", "85÷5="),
    @("29÷5=", "49÷5="),
    @("23÷9=", "44÷5="),
    @("38÷3=", "87÷8="),
    @("75÷8=", "88÷3="),
    @("58÷6=", "14÷5="),
    @("84÷8=", "36÷2="),
    @("36÷7=", "70÷9="),
    @("49÷3=", "50÷6="),
    @("95÷6=", "98÷2="),
    @("23÷5=", "81÷3="),
    @("25÷9=", "14÷9="),
    @("19÷3=", "74÷6="),
    @("56÷4=", "92÷3="),
    @("41÷6=", "85÷9="),
    @("30÷6=", "53÷9="),
    @("42÷6=", "67÷9="),
    @("15÷8=", "81÷6="),
    @("25÷6=", "11÷3="),
    @("86÷8=", "21÷2="),
    @("91÷4=", "43÷6="),
    @("79÷8=", "16÷5="),
    @("49÷8=", "34÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
